$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 72
$ws.Range("F2").Value = 53
$ws.Range("H2").Value = 53

# Row 10
$ws.Range("E10").Value = 269
$ws.Range("F10").Value = 135
$ws.Range("H10").Value = 135

# Row 11
$ws.Range("E11").Value = 194

# Row 12
$ws.Range("E12").Value = 283
$ws.Range("F12").Value = 170
$ws.Range("H12").Value = 170

# Row 13
$ws.Range("E13").Value = 91
$ws.Range("F13").Value = 45
$ws.Range("H13").Value = 45

# Row 14
$ws.Range("E14").Value = 80

# Row 15
$ws.Range("E15").Value = 108

# Row 16
$ws.Range("E16").Value = 114

# Row 20
$ws.Range("E20").Value = 63

# Row 21
$ws.Range("E21").Value = 86

# Row 24
$ws.Range("E24").Value = 131

# Row 25
$ws.Range("E25").Value = 138
$ws.Range("F25").Value = 68
$ws.Range("H25").Value = 68

# Row 26
$ws.Range("E26").Value = 79

# Row 27
$ws.Range("E27").Value = 184

# Row 28
$ws.Range("E28").Value = 113
$ws.Range("F28").Value = 39
$ws.Range("H28").Value = 39

# Row 30
$ws.Range("E30").Value = 133

# Row 32
$ws.Range("E32").Value = 120

# Row 33
$ws.Range("E33").Value = 172

# Row 34
$ws.Range("E34").Value = 130

# Row 36
$ws.Range("E36").Value = 35
$ws.Range("F36").Value = 25
$ws.Range("H36").Value = 25

# Row 37
$ws.Range("E37").Value = 90
$ws.Range("F37").Value = 53
$ws.Range("H37").Value = 53

# Row 40
$ws.Range("E40").Value = 167
$ws.Range("F40").Value = 82
$ws.Range("H40").Value = 82

# Row 41
$ws.Range("E41").Value = 231

# Row 46
$ws.Range("E46").Value = 169
$ws.Range("F46").Value = 105
$ws.Range("H46").Value = 105

# Row 47
$ws.Range("E47").Value = 269

# Row 48
$ws.Range("E48").Value = 127

# Row 49
$ws.Range("E49").Value = 151
$ws.Range("F49").Value = 74
$ws.Range("H49").Value = 74

# Row 50
$ws.Range("E50").Value = 128

# Row 51
$ws.Range("E51").Value = 125

# Row 52
$ws.Range("E52").Value = 13
